# Update "want to go" counts (column F) for several events on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1872
$ws1.Range("F3").Value  = 489
$ws1.Range("F6").Value  = 2577
$ws1.Range("F8").Value  = 90
$ws1.Range("F10").Value = 1530
$ws1.Range("F11").Value = 529
$ws1.Range("F22").Value = 58
$ws1.Range("F23").Value = 1648
$ws1.Range("F29").Value = 418

# --- Sheet "全部类型" ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value  = 1872
$ws2.Range("F4").Value  = 489
$ws2.Range("F7").Value  = 2577
$ws2.Range("F9").Value  = 90
$ws2.Range("F11").Value = 1530
$ws2.Range("F12").Value = 529
$ws2.Range("F23").Value = 58
$ws2.Range("F24").Value = 1648
$ws2.Range("F30").Value = 418
